$wb = $excel.ActiveWorkbook

# summary sheet (Worksheets.Item(1))
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 2850.0000285
$ws.Range("D2").Value = 60998.39019467979

# bus sheet (Worksheets.Item(2))
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = -7.014224642325977
$ws.Range("B3").Value = -7.103140916682147
$ws.Range("B4").Value = -6.404042958386182
$ws.Range("B5").Value = -10.23923767718209
$ws.Range("B6").Value = -10.39735582046898
$ws.Range("B7").Value = -13.14961201638763
$ws.Range("B8").Value = -13.65044939529142
$ws.Range("B9").Value = -15.27657009612768
$ws.Range("B10").Value = -8.409271532435426
$ws.Range("B11").Value = -10.34058346804791
$ws.Range("B12").Value = -2.790266574471174
$ws.Range("B13").Value = -2.059230895068654
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 1.32057356624658
$ws.Range("B16").Value = 10.28199939289767
$ws.Range("B17").Value = 9.470095196720891
$ws.Range("B18").Value = 14.22823343942049
$ws.Range("B19").Value = 15.71124271994595
$ws.Range("B20").Value = 8.253737254026634
$ws.Range("B21").Value = 9.264510898442829
$ws.Range("B22").Value = 16.54779436894681
$ws.Range("B23").Value = 22.72859127851293
$ws.Range("B24").Value = 10.60789884112406
$ws.Range("B25").Value = 3.904933133453839

# branch sheet (Worksheets.Item(6))
$ws = $wb.Worksheets.Item(6)
$ws.Range("D2").Value = -228.8645141549848
$ws.Range("D3").Value = 54.9639733419143
$ws.Range("D4").Value = -179.7457971026728
$ws.Range("D5").Value = -159.1080478896982
$ws.Range("D6").Value = -365.6452054346382
$ws.Range("D7").Value = 46.22338805270467
$ws.Range("D8").Value = 29.40995031280614
$ws.Range("D9").Value = -320.6377552126719
$ws.Range("D10").Value = -30.79935817343193
$ws.Range("D11").Value = 214.4524139594521
$ws.Range("D12").Value = -171.6452044945984
$ws.Range("D13").Value = -102.309535152862
$ws.Range("D14").Value = -214.0378745616294
$ws.Range("D15").Value = -108.5488089415843
$ws.Range("D16").Value = 81.90983487886993
$ws.Range("D17").Value = -140.8919581099992
$ws.Range("D18").Value = -5.04246184664594
$ws.Range("D19").Value = -56.3728982163642
$ws.Range("D20").Value = -72.59659034151359
$ws.Range("D21").Value = -223.1811277284571
$ws.Range("D22").Value = -81.03602801808567
$ws.Range("D23").Value = 43.20064256656804
$ws.Range("D24").Value = 11.16461687857928
$ws.Range("D25").Value = -44.54880830158424
$ws.Range("D26").Value = -108.5488089415842
$ws.Range("D27").Value = 69.87784388796968
$ws.Range("D28").Value = 91.9023852068315
$ws.Range("D29").Value = -52.18002331578169
$ws.Range("D30").Value = -75.50495625559178
$ws.Range("D31").Value = -1.122156822030282
$ws.Range("D32").Value = -44.54880830158424
$ws.Range("D33").Value = -223.1811277284571
$ws.Range("D34").Value = -56.37289821636418

# transformer sheet (Worksheets.Item(7))
$ws = $wb.Worksheets.Item(7)
$ws.Range("D2").Value = -116.8893172852906
$ws.Range("D3").Value = -132.0966826668487
$ws.Range("D4").Value = -157.0654223621697
$ws.Range("D5").Value = -214.4524139594521
$ws.Range("D6").Value = -172.2727877437279
